$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Dashboard Page "
$ws.Range("B10").Value = "TC5"
$ws.Range("C10").Value = "Verify user profile details name on all section "
$ws.Range("D10").Value = "Medium"
$ws.Range("E10").Value = "Sanity"

[void]$ws.Range("C10").Select()
